$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 3357.8462
$ws.Range("I5").Value = 3758.4
$ws.Range("J5").Value = 2022.6666
$ws.Range("K5").Value = 3758.4
$ws.Range("L5").Value = 2022.6666
$ws.Range("M5").Value = -3643.4
$ws.Range("N5").Value = -2252.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2833.6
$ws.Range("I100").Value = 2295.25
$ws.Range("K100").Value = 2295.25
$ws.Range("M100").Value = -1754.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1924.8182
$ws.Range("I132").Value = 1923.4259
$ws.Range("K132").Value = 5770.2777
$ws.Range("M132").Value = -3240.2777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 183.46153
$ws.Range("I4").Value = 167.14285
$ws.Range("J4").Value = 202.5
$ws.Range("K4").Value = 167.14285
$ws.Range("L4").Value = 202.5
$ws.Range("M4").Value = -51.14285000000001
$ws.Range("N4").Value = -434.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13484.958
$ws.Range("I32").Value = 9697.316999999999
$ws.Range("J32").Value = 35669.715
$ws.Range("K32").Value = 9697.316999999999
$ws.Range("L32").Value = 35669.715
$ws.Range("M32").Value = -9410.316999999999
$ws.Range("N32").Value = -36243.715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 5250
$ws.Range("J88").Value = 6000
$ws.Range("L88").Value = 6000
$ws.Range("N88").Value = -6812

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 5250
$ws.Range("J91").Value = 6000
$ws.Range("L91").Value = 6000
$ws.Range("N91").Value = -8808

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 807.129
$ws.Range("I97").Value = 715.03845
$ws.Range("J97").Value = 1286
$ws.Range("K97").Value = 715.03845
$ws.Range("L97").Value = 1286
$ws.Range("M97").Value = -219.03845
$ws.Range("N97").Value = -2278

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 698.4286
$ws.Range("I22").Value = 572.73334
$ws.Range("K22").Value = 572.73334
$ws.Range("M22").Value = -399.73334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3651.0938
$ws.Range("I86").Value = 1440.8572
$ws.Range("J86").Value = 7870.636
$ws.Range("K86").Value = 1440.8572
$ws.Range("L86").Value = 7870.636
$ws.Range("M86").Value = -317.8571999999999
$ws.Range("N86").Value = -10116.636

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 3651.0938
$ws.Range("I89").Value = 1440.8572
$ws.Range("J89").Value = 7870.636
$ws.Range("K89").Value = 7204.286
$ws.Range("L89").Value = 39353.18
$ws.Range("M89").Value = -1588.286
$ws.Range("N89").Value = -50585.18

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1073.2258
$ws.Range("I94").Value = 523.55
$ws.Range("K94").Value = 523.55
$ws.Range("M94").Value = -72.54999999999995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3733.9038
$ws.Range("I134").Value = 2839.2559
$ws.Range("J134").Value = 8008.3335
$ws.Range("K134").Value = 8517.7677
$ws.Range("L134").Value = 24025.0005
$ws.Range("M134").Value = -5982.7677
$ws.Range("N134").Value = -29095.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11089.896
$ws.Range("J31").Value = 16200.1875
$ws.Range("L31").Value = 16200.1875
$ws.Range("N31").Value = -16790.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 11089.896
$ws.Range("J34").Value = 16200.1875
$ws.Range("L34").Value = 16200.1875
$ws.Range("N34").Value = -16604.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3917.7334
$ws.Range("I58").Value = 3078.1538
$ws.Range("J58").Value = 9375
$ws.Range("K58").Value = 3078.1538
$ws.Range("L58").Value = 9375
$ws.Range("M58").Value = -2875.1538
$ws.Range("N58").Value = -9781

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3917.7334
$ws.Range("I136").Value = 3078.1538
$ws.Range("J136").Value = 9375
$ws.Range("K136").Value = 9234.4614
$ws.Range("L136").Value = 28125
$ws.Range("M136").Value = -6684.4614
$ws.Range("N136").Value = -33225

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1576.125
$ws.Range("I5").Value = 929.75
$ws.Range("K5").Value = 2789.25
$ws.Range("M5").Value = -2677.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1596.375
$ws.Range("J113").Value = 1653.1428
$ws.Range("L113").Value = 4959.428400000001
$ws.Range("N113").Value = -9299.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2810.5
$ws.Range("I132").Value = 3048.4
$ws.Range("J132").Value = 2414
$ws.Range("K132").Value = 27435.6
$ws.Range("L132").Value = 21726
$ws.Range("M132").Value = -24905.6
$ws.Range("N132").Value = -26786

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1576.125
$ws.Range("I135").Value = 929.75
$ws.Range("K135").Value = 8367.75
$ws.Range("M135").Value = -5832.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1772.762
$ws.Range("I140").Value = 1342.6666
$ws.Range("J140").Value = 1844.4445
$ws.Range("K140").Value = 4027.9998
$ws.Range("L140").Value = 5533.333500000001
$ws.Range("M140").Value = 1152.0002
$ws.Range("N140").Value = -15893.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 2401.7273
$ws.Range("J2").Value = 444.83334
$ws.Range("L2").Value = 444.83334
$ws.Range("N2").Value = -670.83334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2949.6667
$ws.Range("I80").Value = 2899.5
$ws.Range("J80").Value = 2974.75
$ws.Range("K80").Value = 2899.5
$ws.Range("L80").Value = 2974.75
$ws.Range("M80").Value = -1901.5
$ws.Range("N80").Value = -4970.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2949.6667
$ws.Range("I83").Value = 2899.5
$ws.Range("J83").Value = 2974.75
$ws.Range("K83").Value = 14497.5
$ws.Range("L83").Value = 14873.75
$ws.Range("M83").Value = -9505.5
$ws.Range("N83").Value = -24857.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2538.5217
$ws.Range("I113").Value = 2914.8333
$ws.Range("K113").Value = 2914.8333
$ws.Range("M113").Value = -744.8332999999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2941.41
$ws.Range("I122").Value = 1998.2106
$ws.Range("K122").Value = 5994.6318
$ws.Range("M122").Value = -3544.6318

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5038.351
$ws.Range("I132").Value = 4191.054
$ws.Range("K132").Value = 12573.162
$ws.Range("M132").Value = -10043.162

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 24573.166
$ws.Range("J136").Value = 24573.166
$ws.Range("L136").Value = 73719.49800000001
$ws.Range("N136").Value = -78819.49800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3532.8
$ws.Range("I22").Value = 4432.7144
$ws.Range("J22").Value = 1433
$ws.Range("K22").Value = 4432.7144
$ws.Range("L22").Value = 1433
$ws.Range("M22").Value = -4137.7144
$ws.Range("N22").Value = -2023

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3532.8
$ws.Range("I27").Value = 4432.7144
$ws.Range("J27").Value = 1433
$ws.Range("K27").Value = 4432.7144
$ws.Range("L27").Value = 1433
$ws.Range("M27").Value = -4325.7144
$ws.Range("N27").Value = -1647

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 401.52
$ws.Range("J55").Value = 565.5
$ws.Range("L55").Value = 565.5
$ws.Range("N55").Value = -911.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2679
$ws.Range("I82").Value = 945.25
$ws.Range("J82").Value = 3669.7144
$ws.Range("K82").Value = 945.25
$ws.Range("L82").Value = 3669.7144
$ws.Range("M82").Value = -584.25
$ws.Range("N82").Value = -4391.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2679
$ws.Range("I85").Value = 945.25
$ws.Range("J85").Value = 3669.7144
$ws.Range("K85").Value = 945.25
$ws.Range("L85").Value = 3669.7144
$ws.Range("M85").Value = 302.75
$ws.Range("N85").Value = -6165.7144

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 4343
$ws.Range("I93").Value = 2667.3333
$ws.Range("K93").Value = 2667.3333
$ws.Range("M93").Value = -1419.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4116.222
$ws.Range("I132").Value = 2206.7778
$ws.Range("J132").Value = 6025.6665
$ws.Range("K132").Value = 6620.3334
$ws.Range("L132").Value = 18076.9995
$ws.Range("M132").Value = -4090.3334
$ws.Range("N132").Value = -23136.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2843.8235
$ws.Range("I122").Value = 2720.5518
$ws.Range("K122").Value = 8161.655400000001
$ws.Range("M122").Value = -5711.655400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4115.278
$ws.Range("I132").Value = 3556.516
$ws.Range("K132").Value = 10669.548
$ws.Range("M132").Value = -8139.548000000001
